$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grading pass: award full marks (10/10) on the two rubric rows that were
# previously left ungraded ("whoPurchasedProduct() method" and
# "findAllBrands()" under the CustomerMapping Class section).
$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Leave the sheet scrolled/selected where the grader finished working.
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E24").Select() | Out-Null
